# The edit collapses the final "ArticleParagraph.4" title+content pair of
# paragraphs into a single empty paragraph (pPr style20, one run with an
# empty rPr and no text) -- mirroring the blank separator paragraphs that
# already sit between the other article sections.

$d = $word.ActiveDocument

$titleText   = "`${Article.ArticleParagraph.4.title}"
$contentText = "`${Article.ArticleParagraph.4.content}"

# Locate the title / content paragraphs by their placeholder text (Range.Text
# includes the trailing paragraph-mark character, so trim it before compare).
$titleParaIdx   = -1
$contentParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq $titleText)   { $titleParaIdx = $i }
    if ($txt -eq $contentText) { $contentParaIdx = $i }
}

# Remove the title paragraph, the content paragraph, and the paragraph mark
# that separates them from the blank paragraph right before the title -- this
# merges all three into the single preceding (already-blank) paragraph mark.
$beforePara = $d.Paragraphs.Item($titleParaIdx - 1)
$lastPara   = $d.Paragraphs.Item($contentParaIdx)
$killRange  = $d.Range($beforePara.Range.End, $lastPara.Range.End)
$killRange.Delete()

# Re-create a fresh empty paragraph after the (still blank/unformatted)
# separator paragraph, giving back the paragraph that used to hold the title
# text -- now with no text and no direct formatting, i.e. <w:r><w:rPr/></w:r>.
$beforePara2 = $d.Paragraphs.Item($titleParaIdx - 1)
$beforePara2.Range.InsertParagraphAfter()
